# Ifood-Original.xlsx - re-categorize several tweets in column B (Categoria)
# and reset the view/selection back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the category classification for the rows called out in the diff.
# (shared-string text values: BOM, IMPOSSIVEL SABER, RUIM, IRRELEVANTE, CUPOM)
$ws.Range("B2").Value   = "BOM"
$ws.Range("B7").Value   = "BOM"
$ws.Range("B27").Value  = "IRRELEVANTE"
$ws.Range("B31").Value  = "IRRELEVANTE"
$ws.Range("B37").Value  = "IRRELEVANTE"
$ws.Range("B50").Value  = "BOM"
$ws.Range("B75").Value  = "IRRELEVANTE"
$ws.Range("B98").Value  = "IRRELEVANTE"
$ws.Range("B100").Value = "RUIM"
$ws.Range("B102").Value = "BOM"
$ws.Range("B111").Value = "BOM"
$ws.Range("B112").Value = "RUIM"
$ws.Range("B162").Value = "BOM"
$ws.Range("B167").Value = "BOM"
$ws.Range("B184").Value = "BOM"
$ws.Range("B188").Value = "RUIM"
$ws.Range("B201").Value = "BOM"
$ws.Range("B208").Value = "RUIM"
$ws.Range("B210").Value = "BOM"
$ws.Range("B211").Value = "BOM"
$ws.Range("B217").Value = "BOM"
$ws.Range("B218").Value = "BOM"
$ws.Range("B228").Value = "BOM"
$ws.Range("B242").Value = "BOM"
$ws.Range("B249").Value = "RUIM"
$ws.Range("B255").Value = "BOM"
$ws.Range("B268").Value = "BOM"
$ws.Range("B276").Value = "BOM"
$ws.Range("B299").Value = "IRRELEVANTE"
$ws.Range("B301").Value = "IRRELEVANTE"
$ws.Range("B302").Value = "IRRELEVANTE"

# Scroll back up and move the selection/active cell to B6 (matching the
# saved view state in the workbook).
$ws.Range("A1").Select()
$ws.Range("B6").Select()
